$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns a value as text (forcing string type even for numeric-looking
# strings like "275.94"), then resets the cell style so no stray number-format
# / quote-prefix style is left behind on the cell.
function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "20.490.05"
Set-TextCell $ws "E2" "  +2.77%  "
Set-TextCell $ws "D3" "1.473.65"
Set-TextCell $ws "E3" "  +4.42%  "
Set-TextCell $ws "E4" "  +0.64%  "
Set-TextCell $ws "D5" "0.9592"
Set-TextCell $ws "E5" "  -4.17%  "
Set-TextCell $ws "D6" "275.94"
Set-TextCell $ws "E6" "  +0.15%  "
Set-TextCell $ws "D7" "0.3648"
Set-TextCell $ws "E7" "  -0.79%  "
Set-TextCell $ws "D8" "0.3064"
Set-TextCell $ws "E8" "  -1.17%  "
Set-TextCell $ws "D9" "39.87"
Set-TextCell $ws "E9" "  +0.02%  "
Set-TextCell $ws "E10" "  +1.42%  "
Set-TextCell $ws "E11" "  +1.83%  "
Set-TextCell $ws "E12" "  +0.07%  "
Set-TextCell $ws "D13" "18.18"
Set-TextCell $ws "E13" "  +3.70%  "
Set-TextCell $ws "D14" "5.454"
Set-TextCell $ws "E14" "  -0.06%  "
Set-TextCell $ws "D15" "6.168"
Set-TextCell $ws "E15" "  +0.07%  "
Set-TextCell $ws "E16" "  +1.34%  "
Set-TextCell $ws "D17" "1.473.44"
Set-TextCell $ws "E17" "  +4.15%  "
Set-TextCell $ws "D18" "0.05896"
Set-TextCell $ws "E18" "  +3.97%  "
Set-TextCell $ws "D19" "0.9696"
Set-TextCell $ws "E19" "  -3.14%  "
Set-TextCell $ws "D20" "68.92"
Set-TextCell $ws "E20" "  -2.59%  "
Set-TextCell $ws "D21" "5.464"
Set-TextCell $ws "E21" "  -2.00%  "
Set-TextCell $ws "D22" "14.41"
Set-TextCell $ws "E22" "  -1.82%  "
Set-TextCell $ws "E23" "  +0.11%  "
Set-TextCell $ws "D24" "2.249"
Set-TextCell $ws "E24" "  +0.74%  "
Set-TextCell $ws "D25" "20.527.94"
Set-TextCell $ws "E25" "  +2.79%  "
Set-TextCell $ws "D26" "141.74"
Set-TextCell $ws "E26" "  +6.71%  "
Set-TextCell $ws "D27" "2.128"
Set-TextCell $ws "E27" "  -6.37%  "
Set-TextCell $ws "E28" "  +0.09%  "
Set-TextCell $ws "D29" "1.630.07"
Set-TextCell $ws "E29" "  +3.49%  "
Set-TextCell $ws "D30" "113.36"
Set-TextCell $ws "E30" "  +3.22%  "
Set-TextCell $ws "D31" "3.885"
Set-TextCell $ws "E31" "  -0.46%  "
Set-TextCell $ws "B32" "ImmutableX"
Set-TextCell $ws "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D32" "0.8114"
Set-TextCell $ws "E32" "  +0.50%  "
Set-TextCell $ws "B33" "Filecoin"
Set-TextCell $ws "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D33" "4.956"
Set-TextCell $ws "E33" "  -4.32%  "
Set-TextCell $ws "D34" "0.07871"
Set-TextCell $ws "E34" "  +1.49%  "
Set-TextCell $ws "D35" "1.255"
Set-TextCell $ws "E35" "  +14.40%  "
Set-TextCell $ws "D36" "1.524"
Set-TextCell $ws "E36" "  +3.91%  "
Set-TextCell $ws "D37" "0.05731"
Set-TextCell $ws "E37" "  -1.12%  "
Set-TextCell $ws "D38" "4.740"
Set-TextCell $ws "E38" "  -2.78%  "
Set-TextCell $ws "B39" "FraxShare"
Set-TextCell $ws "C39" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D39" "7.660"
Set-TextCell $ws "E39" "  -4.84%  "
Set-TextCell $ws "D40" "0.02036"
Set-TextCell $ws "E40" "  -0.21%  "
Set-TextCell $ws "B41" "Frax"
Set-TextCell $ws "C41" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws "D41" "0.9632"
Set-TextCell $ws "E41" "  -3.74%  "
Set-TextCell $ws "D42" "10.41"
Set-TextCell $ws "E42" "  +0.02%  "
Set-TextCell $ws "D43" "0.1874"
Set-TextCell $ws "E43" "  -0.08%  "
Set-TextCell $ws "D44" "0.5275"
Set-TextCell $ws "E44" "  -0.15%  "
Set-TextCell $ws "E45" "  -0.77%  "
Set-TextCell $ws "D46" "12.10"
Set-TextCell $ws "E46" "  -1.33%  "
Set-TextCell $ws "D47" "116.75"
Set-TextCell $ws "E47" "  -0.02%  "
Set-TextCell $ws "D48" "0.5175"
Set-TextCell $ws "E48" "  +0.22%  "
Set-TextCell $ws "D49" "1.767"
Set-TextCell $ws "E49" "  +0.43%  "
Set-TextCell $ws "D50" "0.06451"
Set-TextCell $ws "E50" "  +4.34%  "
Set-TextCell $ws "D51" "0.9888"
Set-TextCell $ws "E51" "  -1.18%  "
